# Updated Excel files with primary track settings.
# Inserts 4 new dataset rows (BF-C2DL-HSC, BF-C2DL-MuSC, DIC-C2DH-HeLa, PhC-C2DH-U373)
# into the alphabetically-sorted settings table and refreshes the shared
# numeric setting values (columns B-J) for every dataset row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final, alphabetically-sorted list of datasets (column A, rows 2..14)
$datasets = @(
    "BF-C2DL-HSC",
    "BF-C2DL-MuSC",
    "DIC-C2DH-HeLa",
    "Fluo-C2DL-MSC",
    "Fluo-C3DH-A549",
    "Fluo-C3DH-H157",
    "Fluo-C3DL-MDA231",
    "Fluo-N2DH-GOWT1",
    "Fluo-N2DL-HeLa",
    "Fluo-N3DH-CE",
    "Fluo-N3DH-CHO",
    "PhC-C2DH-U373",
    "PhC-C2DL-PSC"
)

# Shared setting values (columns B,C,D,E,G,H,I,J) - identical for every row.
# Column F (SegClipping) is the literal number 1.
$settingValues = @("37.4027", "4.6186", "0.93389", "0.0071904", "0.023956", "7.496", "151.1994", "26.6334")

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $datasets.Count; $i++) {
    $r = $i + 2

    Set-TextCell $ws.Cells.Item($r, 1) $datasets[$i]

    Set-TextCell $ws.Cells.Item($r, 2) $settingValues[0]
    Set-TextCell $ws.Cells.Item($r, 3) $settingValues[1]
    Set-TextCell $ws.Cells.Item($r, 4) $settingValues[2]
    Set-TextCell $ws.Cells.Item($r, 5) $settingValues[3]

    $ws.Cells.Item($r, 6).Value = 1

    Set-TextCell $ws.Cells.Item($r, 7) $settingValues[4]
    Set-TextCell $ws.Cells.Item($r, 8) $settingValues[5]
    Set-TextCell $ws.Cells.Item($r, 9) $settingValues[6]
    Set-TextCell $ws.Cells.Item($r, 10) $settingValues[7]
}

Write-Host "Done updating Settings-allGT table."
